# fix: prevent hidden columns from being labeled upon detecting changes
#
# Column L ("Aenderung") was incorrectly flagged with the "AENDERUNG" marker
# (style 7 = bold/gold centered text on grey fill) for a batch of rows whose
# apparent delta only came from a hidden column. Those L cells should go back
# to being blank, using the plain grey "no change" look (style 4) instead.
#
# Two of those rows (212 and 218) are also the first row of a new group and
# need their whole row's look-and-feel (styles 5/7 -> 2/3/4) brought in line
# with the other "new group" rows (e.g. row 2), in addition to clearing
# column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that need a full-row restyle (new-group first rows) -------------
$fullRestyleRows = @(212, 218)

$templateRow = $ws.Range("A2:V2")
foreach ($r in $fullRestyleRows) {
    $destRow = $ws.Range("A" + $r + ":V" + $r)
    $templateRow.Copy()
    $destRow.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# --- Rows that only need column L's mistaken "AENDERUNG" flag removed -----
$simpleRows = @(80,98,99,100,101,102,103,104,105,106,108,111,112,113,114,117,118,119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,134,135,136,139,140,143,149,150,151,158,159,165,166,167,168,169,175,187,188,189,196,223,224)

$lTemplate = $ws.Range("L2")
foreach ($r in $simpleRows) {
    $cell = $ws.Range("L" + $r)
    $lTemplate.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# Clear the stray "AENDERUNG" text out of every affected column-L cell,
# including the two full-restyle rows above.
$allLRows = $fullRestyleRows + $simpleRows
foreach ($r in $allLRows) {
    $ws.Range("L" + $r).ClearContents()
}
